$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Paragraph 1 (style "FirstParagraph") — the "elephant in the room"
# paragraph. The whole paragraph is a single run with no special
# character formatting, so replace its whole text in one shot.
# ------------------------------------------------------------------

$oldPara1 = "I want to talk about the elephant in the room. We the MARC community are grounded in the priciples of maximizing accessing to research careers for URM students. That has been our goal from day one. The events such as the death of george ffloyd, ahmaud arbury, breonna taylor and many other reminds us that there are societal hurdles our students face when they are outside of school. Also, if we listen to the POC in academia we find that their path is similarly laden with racial bias, profiling, a lack of support, and in many cases outright distain."
$newPara1 = "I planned this session to be about your goals, but I want to talk about the elephant in the room. The MARC community is grounded in the principles of maximizing accessing to research careers for URM students. That has been our goal from day one. The death of George Floyd, Ahmaud Arbury, Breonna Taylor and many others reminds us that there are societal hurdles that students of color face when they are outside of this university. Also, if we listen to the POC in academia we find that their path is similarly laden with racial bias, profiling, a lack of support, and in many cases outright distain for the position that they have worked hard to attain."

$found1 = $false
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd("`r`a") -eq $oldPara1) {
        $p.Range.Text = $newPara1
        $found1 = $true
        break
    }
}
if (-not $found1) {
    # Fallback: plain find/replace if the exact whole-paragraph match failed
    # for some reason (e.g. unexpected trailing whitespace).
    $found1 = $d.Content.Find.Execute($oldPara1, $true, $false, $false, $false,
                                       $false, $true, 1, $false, $newPara1, 2)
}
Write-Output "paragraph1 replaced: $found1"

# ------------------------------------------------------------------
# Paragraph 2 (style "BodyText") — the "COVID-19" paragraph. Also a
# single run with no special character formatting.
# ------------------------------------------------------------------

$oldPara2 = "COVID-19 and valid and ongoing protests in support of BLM have created difficult circumstances and conversations. I do not want to shy away from those coversations. I would like you all to feel comfortable in this community sharing your thoughts, concerns, and wants. That is to say this is a place were we respect the thoughts and ideas of other fully."
$newPara2 = "COVID-19 and protests in support of BLM have created spaces where we need to have difficult conversations about our health, our safety, and our rights as people and our rights as researchers in STEM. I do not want to shy away from those coversations. I would like you all to feel comfortable in this community sharing your thoughts, concerns, and wants. That is to say this is a place were we respect the thoughts and ideas of other fully. We can take as much time as we need to discuss because what we are doing for research is changing rapidly, but we can always develop and discuss" + " " + [char]0x201C + "Why" + [char]0x201D + " " + "we do research, which can motivate, activate, and embolden us to push further than we could before.I want you all to know that I and the MARC coordinator team are here to support you whenever you need it, but right now, if you feel comfortable, I" + [char]0x2019 + "d like to hear your thoughts and concerns."

$found2 = $false
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd("`r`a") -eq $oldPara2) {
        $p.Range.Text = $newPara2
        $found2 = $true
        break
    }
}
if (-not $found2) {
    $found2 = $d.Content.Find.Execute($oldPara2, $true, $false, $false, $false,
                                       $false, $true, 1, $false, $newPara2, 2)
}
Write-Output "paragraph2 replaced: $found2"
